$wb = $excel.ActiveWorkbook

# Sheet "展览" (first sheet)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 281
$ws1.Range("F4").Value = 947
$ws1.Range("F6").Value = 52

# Sheet "全部类型" (fourth sheet)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 281
$ws4.Range("F5").Value = 947
$ws4.Range("F7").Value = 52
